# This script updates the "Southwest (NM)_A" team state-transition matrix on Sheet1.
# Each row represents the probability distribution of transitioning from the
# "Starting_State" in column A to each of the possible next states (columns B:S),
# derived from simulated game counts (games added per the "added more games" commit).
# Probabilities are written as simple fractions (observed transitions / total games
# simulated from that starting state) so that every updated row sums to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Starting_State = Af0 (3 simulated games)
$ws.Range("B2").Value = 1/3
$ws.Range("C2").Value = 2/3

# Row 3 - Starting_State = Af1 (4 simulated games)
$ws.Range("P3").Value = 3/4
$ws.Range("S3").Value = 1/4

# Row 6 - Starting_State = Ai0 (10 simulated games)
$ws.Range("B6").Value = 2/10
$ws.Range("F6").Value = 2/10
$ws.Range("J6").Value = 2/10
$ws.Range("O6").Value = 1/10
$ws.Range("R6").Value = 1/10
$ws.Range("S6").Value = 2/10

# Row 7 - Starting_State = Ai1 (5 simulated games)
$ws.Range("B7").Value = 1/5
$ws.Range("Q7").Value = 1/5
$ws.Range("R7").Value = 1/5
$ws.Range("S7").Value = 2/5

# Row 8 - Starting_State = Ai2 (8 simulated games)
$ws.Range("J8").Value = 1/8
$ws.Range("Q8").Value = 3/8
$ws.Range("R8").Value = 4/8

# Row 9 - Starting_State = Ai3 (7 simulated games)
$ws.Range("B9").Value = 1/7
$ws.Range("J9").Value = 2/7
$ws.Range("Q9").Value = 2/7
$ws.Range("S9").Value = 2/7

# Row 10 - Starting_State = Ar0 (55 simulated games)
$ws.Range("B10").Value = 1/55
$ws.Range("D10").Value = 1/55
$ws.Range("F10").Value = 2/55
$ws.Range("J10").Value = 10/55
$ws.Range("O10").Value = 1/55
$ws.Range("Q10").Value = 15/55
$ws.Range("R10").Value = 4/55
$ws.Range("S10").Value = 21/55

# Row 11 - Starting_State = Bf0 (7 simulated games)
$ws.Range("G11").Value = 1/7
$ws.Range("J11").Value = 1/7
$ws.Range("K11").Value = 2/7
$ws.Range("L11").Value = 3/7

# Row 12 - Starting_State = Bf1 (1 simulated game)
$ws.Range("G12").Value = 1

# Row 13 - Starting_State = Bf2 (1 simulated game)
$ws.Range("G13").Value = 1

# Row 15 - Starting_State = Bi0 (9 simulated games)
$ws.Range("F15").Value = 1/9
$ws.Range("I15").Value = 1/9
$ws.Range("J15").Value = 4/9
$ws.Range("O15").Value = 2/9
$ws.Range("S15").Value = 1/9

# Row 16 - Starting_State = Bi1 (3 simulated games)
$ws.Range("H16").Value = 1/3
$ws.Range("K16").Value = 1/3
$ws.Range("S16").Value = 1/3

# Row 17 - Starting_State = Bi2 (20 simulated games)
$ws.Range("F17").Value = 1/20
$ws.Range("H17").Value = 2/20
$ws.Range("I17").Value = 1/20
$ws.Range("J17").Value = 11/20
$ws.Range("O17").Value = 2/20
$ws.Range("S17").Value = 3/20

# Row 18 - Starting_State = Bi3 (10 simulated games)
$ws.Range("F18").Value = 1/10
$ws.Range("I18").Value = 2/10
$ws.Range("J18").Value = 7/10

# Row 19 - Starting_State = Br0 (35 simulated games)
$ws.Range("F19").Value = 2/35
$ws.Range("H19").Value = 5/35
$ws.Range("I19").Value = 3/35
$ws.Range("J19").Value = 17/35
$ws.Range("K19").Value = 3/35
$ws.Range("M19").Value = 1/35
$ws.Range("O19").Value = 2/35
$ws.Range("S19").Value = 2/35

